# HOTFIX change date format in example
# Replace the example "date" columns (K:P, rows 2-10) with ISO-formatted
# (yyyy-mm-dd) sample dates instead of the old dd.mm.yyyy ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("2021-01-20", "2021-01-21", "2021-01-22", "2021-01-23", "2021-01-24", "2021-01-25")
$cols = @("K", "L", "M", "N", "O", "P")

for ($row = 2; $row -le 10; $row++) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $dates[$i]
    }
}
